$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A4").Value = "Nycthemeral timing"
$ws.Range("A5").Value = "Nycthemeral regularity"
